$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Style = 'Normal'
$ws.Range('D2').Value = '37.029.35'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.94%  '
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').Style = 'Normal'
$ws.Range('D3').Value = '2.009.02'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.88%  '
$ws.Range('E3').Style = 'Normal'

$ws.Range('E4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').Style = 'Normal'
$ws.Range('D5').Value = '''225.78'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.19%  '
$ws.Range('E5').Style = 'Normal'

$ws.Range('D6').Style = 'Normal'
$ws.Range('D6').Value = '''0.604'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.94%  '
$ws.Range('E6').Style = 'Normal'

$ws.Range('E7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E7').Style = 'Normal'

$ws.Range('D8').Style = 'Normal'
$ws.Range('D8').Value = '''55.17'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.17%  '
$ws.Range('E8').Style = 'Normal'

$ws.Range('D9').Style = 'Normal'
$ws.Range('D9').Value = '''0.373'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.00%  '
$ws.Range('E9').Style = 'Normal'

$ws.Range('D10').Style = 'Normal'
$ws.Range('D10').Value = '''0.0776'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.05%  '
$ws.Range('E10').Style = 'Normal'

$ws.Range('E11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.93%  '
$ws.Range('E11').Style = 'Normal'

$ws.Range('D12').Style = 'Normal'
$ws.Range('D12').Value = '2.306.21'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.95%  '
$ws.Range('E12').Style = 'Normal'

$ws.Range('D13').Style = 'Normal'
$ws.Range('D13').Value = '''14.03'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.89%  '
$ws.Range('E13').Style = 'Normal'

$ws.Range('D14').Style = 'Normal'
$ws.Range('D14').Value = '''19.72'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.79%  '
$ws.Range('E14').Style = 'Normal'

$ws.Range('D15').Style = 'Normal'
$ws.Range('D15').Value = '''0.733'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.47%  '
$ws.Range('E15').Style = 'Normal'

$ws.Range('D16').Style = 'Normal'
$ws.Range('D16').Value = '''5.14'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.84%  '
$ws.Range('E16').Style = 'Normal'

$ws.Range('D17').Style = 'Normal'
$ws.Range('D17').Value = '2.012.41'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.80%  '
$ws.Range('E17').Style = 'Normal'

$ws.Range('D18').Style = 'Normal'
$ws.Range('D18').Value = '36.955.38'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.86%  '
$ws.Range('E18').Style = 'Normal'

$ws.Range('D19').Style = 'Normal'
$ws.Range('D19').Value = '''6.16'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.81%  '
$ws.Range('E19').Style = 'Normal'

$ws.Range('D20').Style = 'Normal'
$ws.Range('D20').Value = '''68.21'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.18%  '
$ws.Range('E20').Style = 'Normal'

$ws.Range('D21').Style = 'Normal'
$ws.Range('D21').Value = '0.0₃0810'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.19%  '
$ws.Range('E21').Style = 'Normal'

$ws.Range('D22').Style = 'Normal'
$ws.Range('D22').Value = '''222.77'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.25%  '
$ws.Range('E22').Style = 'Normal'

$ws.Range('D23').Style = 'Normal'
$ws.Range('D23').Value = '''0.999'
$ws.Range('D23').Style = 'Normal'

$ws.Range('E24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.77%  '
$ws.Range('E24').Style = 'Normal'

$ws.Range('E25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.72%  '
$ws.Range('E25').Style = 'Normal'

$ws.Range('D26').Style = 'Normal'
$ws.Range('D26').Value = '''164.54'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.33%  '
$ws.Range('E26').Style = 'Normal'

$ws.Range('D27').Style = 'Normal'
$ws.Range('D27').Value = '''8.93'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.16%  '
$ws.Range('E27').Style = 'Normal'

$ws.Range('E28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.32%  '
$ws.Range('E28').Style = 'Normal'

$ws.Range('D29').Style = 'Normal'
$ws.Range('D29').Value = '''18.54'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.73%  '
$ws.Range('E29').Style = 'Normal'

$ws.Range('E30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.36%  '
$ws.Range('E30').Style = 'Normal'

$ws.Range('D31').Style = 'Normal'
$ws.Range('D31').Value = '''0.116'
$ws.Range('D31').Style = 'Normal'

$ws.Range('D32').Style = 'Normal'
$ws.Range('D32').Value = '''4.38'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.17%  '
$ws.Range('E32').Style = 'Normal'

$ws.Range('D33').Style = 'Normal'
$ws.Range('D33').Value = '''0.0599'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.63%  '
$ws.Range('E33').Style = 'Normal'

$ws.Range('D34').Style = 'Normal'
$ws.Range('D34').Value = '''4.43'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.21%  '
$ws.Range('E34').Style = 'Normal'

$ws.Range('D35').Style = 'Normal'
$ws.Range('D35').Value = '''2.32'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.91%  '
$ws.Range('E35').Style = 'Normal'

$ws.Range('E36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.04%  '
$ws.Range('E36').Style = 'Normal'

$ws.Range('E37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.47%  '
$ws.Range('E37').Style = 'Normal'

$ws.Range('E38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.86%  '
$ws.Range('E38').Style = 'Normal'

$ws.Range('E39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.06%  '
$ws.Range('E39').Style = 'Normal'

$ws.Range('E40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.67%  '
$ws.Range('E40').Style = 'Normal'

$ws.Range('D41').Style = 'Normal'
$ws.Range('D41').Value = '''0.0211'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.95%  '
$ws.Range('E41').Style = 'Normal'

$ws.Range('B42').Style = 'Normal'
$ws.Range('B42').Value = 'Aave'
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').Style = 'Normal'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D42').Value = '''94.34'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.53%  '
$ws.Range('E42').Style = 'Normal'

$ws.Range('B43').Style = 'Normal'
$ws.Range('B43').Value = 'FTXToken'
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').Style = 'Normal'
$ws.Range('C43').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D43').Value = '''4.26'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Style = 'Normal'
$ws.Range('E43').Value = '  +20.12%  '
$ws.Range('E43').Style = 'Normal'

$ws.Range('E44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.19%  '
$ws.Range('E44').Style = 'Normal'

$ws.Range('E45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.26%  '
$ws.Range('E45').Style = 'Normal'

$ws.Range('B46').Style = 'Normal'
$ws.Range('B46').Value = 'TrustWalletToken'
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').Style = 'Normal'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D46').Value = '''1.12'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.71%  '
$ws.Range('E46').Style = 'Normal'

$ws.Range('B47').Style = 'Normal'
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').Style = 'Normal'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D47').Value = '''15.85'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.40%  '
$ws.Range('E47').Style = 'Normal'

$ws.Range('D48').Style = 'Normal'
$ws.Range('D48').Value = '''0.995'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.77%  '
$ws.Range('E48').Style = 'Normal'

$ws.Range('D49').Style = 'Normal'
$ws.Range('D49').Value = '''7.04'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.18%  '
$ws.Range('E49').Style = 'Normal'

$ws.Range('D50').Style = 'Normal'
$ws.Range('D50').Value = '''2.88'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.22%  '
$ws.Range('E50').Style = 'Normal'

$ws.Range('D51').Style = 'Normal'
$ws.Range('D51').Value = '2.195.69'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.89%  '
$ws.Range('E51').Style = 'Normal'
